# Apply crypto price/volume updates for Sat Jun 15 11:49:13 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.332.35'
$ws.Range("D3").Value = '3.539.09'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.56'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.32'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.47%  '
$ws.Range("D7").Value = '3.536.90'
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("E11").Value = '  -4.19%  '
$ws.Range("D13").Value = '4.137.82'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("E14").Value = '  -4.52%  '
$ws.Range("E15").Value = '  -5.07%  '
$ws.Range("D16").Value = '3.537.40'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").Value = '66.409.21'
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.94'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.92'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '426.51'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.62%  '
$ws.Range("E23").Value = '  -1.39%  '
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("D25").Value = '3.678.05'
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.25'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.48%  '
$ws.Range("E30").Value = '  -1.85%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.49'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.21%  '
$ws.Range("E33").Value = '  -3.95%  '
$ws.Range("E34").Value = '  -1.26%  '
$ws.Range("D35").Value = '3.529.01'
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("E37").Value = '  -3.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.82'
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.63'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.94%  '
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '173.42'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.91'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -7.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.67'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.20%  '
$ws.Range("E47").Value = '  -2.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.10'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -7.44%  '
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("E50").Value = '  -4.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.944'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.45%  '
